$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '64.401.19'
$ws.Range("E2").Value = '  +0.54%  '

# Row 3
$ws.Range("D3").Value = '3.138.97'
$ws.Range("E3").Value = '  +0.16%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '603.02'
$ws.Range("E5").Value = '  -0.95%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.61'
$ws.Range("E6").Value = '  -1.67%  '

# Row 7
$ws.Range("E7").Value = '  +0.00%  '

# Row 8
$ws.Range("D8").Value = '3.135.31'
$ws.Range("E8").Value = '  +0.21%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.524'
$ws.Range("E9").Value = '  +0.63%  '

# Row 10
$ws.Range("E10").Value = '  -0.43%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.38'
$ws.Range("E11").Value = '  +2.75%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.470'
$ws.Range("E12").Value = '  -0.46%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000255'
$ws.Range("E13").Value = '  +1.11%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.23'
$ws.Range("E14").Value = '  -0.27%  '

# Row 15
$ws.Range("D15").Value = '3.653.09'
$ws.Range("E15").Value = '  +0.40%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.119'
$ws.Range("E16").Value = '  +2.47%  '

# Row 17
$ws.Range("D17").Value = '64.369.93'
$ws.Range("E17").Value = '  +0.53%  '

# Row 18
$ws.Range("D18").Value = '3.134.44'
$ws.Range("E18").Value = '  +0.17%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.91'
$ws.Range("E19").Value = '  +0.99%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '482.20'
$ws.Range("E20").Value = '  +0.69%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.61'
$ws.Range("E21").Value = '  -0.63%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.710'
$ws.Range("E22").Value = '  +0.33%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.68'
$ws.Range("E23").Value = '  -0.88%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '86.83'
$ws.Range("E24").Value = '  +3.44%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.37'
$ws.Range("E25").Value = '  -1.56%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.08%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.76'
$ws.Range("E27").Value = '  -1.81%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.36'
$ws.Range("E28").Value = '  -1.31%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.16'
$ws.Range("E29").Value = '  +5.68%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.06'
$ws.Range("E30").Value = '  -2.81%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.112'
$ws.Range("E31").Value = '  -1.01%  '

# Row 32
$ws.Range("E32").Value = '  -0.04%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '26.85'

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.67'
$ws.Range("E34").Value = '  -2.03%  '

# Row 35
$ws.Range("E35").Value = '  -1.49%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.00'
$ws.Range("E36").Value = '  +1.02%  '

# Row 37
$ws.Range("D37").Value = '0.0₃0754'
$ws.Range("E37").Value = '  -0.43%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '52.56'
$ws.Range("E38").Value = '  -0.91%  '

# Row 39
$ws.Range("E39").Value = '  +1.36%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '446.24'
$ws.Range("E40").Value = '  -2.91%  '

# Row 41
$ws.Range("E41").Value = '  -0.02%  '

# Row 42
$ws.Range("E42").Value = '  +0.81%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.24'
$ws.Range("E43").Value = '  -1.52%  '

# Row 44
$ws.Range("D44").Value = '2.868.10'
$ws.Range("E44").Value = '  +0.42%  '

# Row 45
$ws.Range("E45").Value = '  -2.32%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.23'
$ws.Range("E46").Value = '  -2.03%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.45'
$ws.Range("E47").Value = '  -0.02%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.998'
$ws.Range("E48").Value = '  -0.05%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '26.06'
$ws.Range("E49").Value = '  -0.46%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.114'
$ws.Range("E50").Value = '  +0.15%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '121.16'
$ws.Range("E51").Value = '  +2.12%  '
